# Adding the RES Hourly Production Forecast to the Portfolio
#
# The Lookup-date prefix (column D, and the matching shared-string text)
# moves from 29.08.2024 -> 24.09.2024, the Interval timestamps (column A)
# shift forward by the same 26 days, and the Prediction column (C) picks
# up the refreshed forecast values for rows 30-81.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Prediction (column C) values for the rows whose forecast changed.
$cChanges = @{
    30 = 0
    31 = 0
    32 = 0.02
    33 = 0.055
    34 = 0.112
    35 = 0.183
    36 = 0.252
    37 = 0.317
    38 = 0.38
    39 = 0.448
    40 = 0.513
    41 = 0.569
    42 = 0.62
    43 = 0.666
    44 = 0.707
    45 = 0.74
    46 = 0.771
    47 = 0.802
    48 = 0.832
    49 = 0.858
    50 = 0.866
    51 = 0.872
    52 = 0.887
    53 = 0.896
    54 = 0.894
    55 = 0.887
    56 = 0.5659999999999999
    58 = 0.654
    59 = 0.6830000000000001
    60 = 0.677
    61 = 0.661
    62 = 0.65
    63 = 0.623
    64 = 0.611
    65 = 0.581
    66 = 0.515
    67 = 0.468
    68 = 0.425
    69 = 0.378
    70 = 0.321
    71 = 0.27
    72 = 0.211
    73 = 0.159
    74 = 0.119
    75 = 0.08699999999999999
    76 = 0.061
    77 = 0.044
    78 = 0.03
    79 = 0
    80 = 0
    81 = 0
}

# Data rows run from 2 to 96. Column B holds the Interval number, which is
# reused verbatim as the numeric suffix of the Lookup string in column D.
for ($r = 2; $r -le 96; $r++) {
    $intervalNumber = $ws.Cells.Item($r, 2).Value2

    # Column A: shift the 15-minute timestamp forward by 26 days
    # (29-Aug-2024 -> 24-Sep-2024).
    $ws.Cells.Item($r, 1).Value2 = $ws.Cells.Item($r, 1).Value2 + 26

    # Column C: refresh the forecast value where it changed.
    if ($cChanges.ContainsKey($r)) {
        $ws.Cells.Item($r, 3).Value2 = $cChanges[$r]
    }

    # Column D: rebuild the Lookup text with the new date prefix.
    $ws.Cells.Item($r, 4).Value = "24.09.2024" + $intervalNumber
}
